$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("Input")
$wsOutput = $wb.Worksheets.Item("Output")

# Pre-seed formatting for the rows whose content/role is about to change, by
# copying the number/text format from a row that already carries the target
# look (text rows use the "B3" style, the date row uses the "B7" style).
$wsInput.Range("B3").Copy()
$wsInput.Range("B1:B2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsInput.Range("B7").Copy()
$wsInput.Range("B4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Drop the "office" row and shift firstname/middlename/lastname/dateofbirth
# up into place with the new client's data.
$wsInput.Range("A1").Value = "firstname"
$wsInput.Range("B1").Value = "Jhon"
$wsInput.Range("A2").Value = "middlename"
$wsInput.Range("B2").Value = "N"
$wsInput.Range("A3").Value = "lastname"
$wsInput.Range("B3").Value = "Deer"
$wsInput.Range("A4").Value = "dateofbirth"
$wsInput.Range("B4").Value = 36892
$wsInput.Rows.Item(5).Delete()

# Refresh the concatenated name on the Output sheet.
$wsOutput.Range("B1").Value = "Jhon N Deer"

# Input becomes the active sheet/tab; restore the recorded selections.
$wsOutput.Range("C14").Select() | Out-Null
$wsInput.Activate()
$wsInput.Range("E14").Select() | Out-Null
